# Weekly update: insert two new rows of "Brócoli" price data (week of 45021)
# at the top of the Macroferia Regional de Talca block (previously starting
# at row 485), pushing the rest of that block (old rows 485-516) down by two
# rows (new rows 487-518).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at row 485, shifting existing data down.
$ws.Rows.Item(485).Insert()
$ws.Rows.Item(485).Insert()

# --- New row 485 ---
$ws.Cells.Item(485, 1).Value = 5
$ws.Cells.Item(485, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(485, 3).Value = "Maule"
$ws.Cells.Item(485, 4).Value = 45021
$ws.Cells.Item(485, 5).Value = 7
$ws.Cells.Item(485, 6).Value = 100112023
$ws.Cells.Item(485, 7).Value = "Brócoli"
$ws.Cells.Item(485, 8).Value = "Sin especificar"
$ws.Cells.Item(485, 9).Value = "Primera"
$ws.Cells.Item(485, 10).Value = 3000
$ws.Cells.Item(485, 11).Value = 700
$ws.Cells.Item(485, 12).Value = 700
$ws.Cells.Item(485, 13).Value = 700
$ws.Cells.Item(485, 14).Value = "$/unidad"
$ws.Cells.Item(485, 15).Value = "Región del Maule"
$ws.Cells.Item(485, 16).Value = 700
$ws.Cells.Item(485, 17).Value = 1
$ws.Cells.Item(485, 18).Value = "Hortaliza"

# --- New row 486 ---
$ws.Cells.Item(486, 1).Value = 5
$ws.Cells.Item(486, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(486, 3).Value = "Maule"
$ws.Cells.Item(486, 4).Value = 45021
$ws.Cells.Item(486, 5).Value = 7
$ws.Cells.Item(486, 6).Value = 100112023
$ws.Cells.Item(486, 7).Value = "Brócoli"
$ws.Cells.Item(486, 8).Value = "Sin especificar"
$ws.Cells.Item(486, 9).Value = "Segunda"
$ws.Cells.Item(486, 10).Value = 2000
$ws.Cells.Item(486, 11).Value = 600
$ws.Cells.Item(486, 12).Value = 600
$ws.Cells.Item(486, 13).Value = 600
$ws.Cells.Item(486, 14).Value = "$/unidad"
$ws.Cells.Item(486, 15).Value = "Región del Maule"
$ws.Cells.Item(486, 16).Value = 600
$ws.Cells.Item(486, 17).Value = 1
$ws.Cells.Item(486, 18).Value = "Hortaliza"
